$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 597.6
$ws.Range("I2").Value = 624.75
$ws.Range("K2").Value = 624.75
$ws.Range("M2").Value = -511.75
$ws.Range("H6").Value = 114.166664
$ws.Range("I6").Value = 128
$ws.Range("J6").Value = 45
$ws.Range("K6").Value = 384
$ws.Range("L6").Value = 135
$ws.Range("M6").Value = -272
$ws.Range("N6").Value = -359
$ws.Range("H17").Value = 3290.2856
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H19").Value = 109.42857
$ws.Range("I19").Value = 100.5
$ws.Range("J19").Value = 113
$ws.Range("K19").Value = 100.5
$ws.Range("L19").Value = 113
$ws.Range("M19").Value = 74.5
$ws.Range("N19").Value = -463
$ws.Range("H43").Value = 4920.75
$ws.Range("J43").Value = 5338
$ws.Range("L43").Value = 5338
$ws.Range("N43").Value = -5476
$ws.Range("H88").Value = 2893.6
$ws.Range("I88").Value = 980
$ws.Range("J88").Value = 3372
$ws.Range("K88").Value = 980
$ws.Range("L88").Value = 3372
$ws.Range("M88").Value = -574
$ws.Range("N88").Value = -4184
$ws.Range("H91").Value = 2893.6
$ws.Range("I91").Value = 980
$ws.Range("J91").Value = 3372
$ws.Range("K91").Value = 980
$ws.Range("L91").Value = 3372
$ws.Range("M91").Value = 424
$ws.Range("N91").Value = -6180
$ws.Range("H106").Value = 32440.375
$ws.Range("J106").Value = 7500
$ws.Range("L106").Value = 7500
$ws.Range("N106").Value = -8762
$ws.Range("H116").Value = 8943.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 8943.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 8943.5
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -15827.5
$ws.Range("H125").Value = 809.61536
$ws.Range("J125").Value = 687.5
$ws.Range("L125").Value = 6187.5
$ws.Range("N125").Value = -11107.5
$ws.Range("H131").Value = 6783.077
$ws.Range("I131").Value = 1637.6
$ws.Range("J131").Value = 9999
$ws.Range("K131").Value = 4912.799999999999
$ws.Range("L131").Value = 29997
$ws.Range("M131").Value = 127.2000000000007
$ws.Range("N131").Value = -40077
$ws.Range("H132").Value = 1571.9062
$ws.Range("I132").Value = 1548.5555
$ws.Range("K132").Value = 4645.666499999999
$ws.Range("M132").Value = -2115.666499999999
$ws.Range("H138").Value = 12667.071
$ws.Range("J138").Value = 13173.111
$ws.Range("L138").Value = 39519.333
$ws.Range("N138").Value = -49799.333

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2500
$ws.Range("I45").Value = 2500
$ws.Range("K45").Value = 2500
$ws.Range("M45").Value = -2123
$ws.Range("H61").Value = 1328.15
$ws.Range("I61").Value = 1281.2778
$ws.Range("K61").Value = 1281.2778
$ws.Range("M61").Value = -1069.2778
$ws.Range("H63").Value = 6162.647
$ws.Range("I63").Value = 5001.6665
$ws.Range("K63").Value = 5001.6665
$ws.Range("M63").Value = -4315.6665
$ws.Range("H66").Value = 6162.647
$ws.Range("I66").Value = 5001.6665
$ws.Range("K66").Value = 25008.3325
$ws.Range("M66").Value = -21576.3325
$ws.Range("H102").Value = 1913.2
$ws.Range("I102").Value = 1704.25
$ws.Range("K102").Value = 1704.25
$ws.Range("M102").Value = -82.25
$ws.Range("H130").Value = 25000
$ws.Range("J130").Value = 25000
$ws.Range("L130").Value = 25000
$ws.Range("N130").Value = -35040
$ws.Range("H132").Value = 1787.566
$ws.Range("I132").Value = 1745.0193
$ws.Range("K132").Value = 5235.0579
$ws.Range("M132").Value = -2705.0579
$ws.Range("H136").Value = 1328.15
$ws.Range("I136").Value = 1281.2778
$ws.Range("K136").Value = 3843.8334
$ws.Range("M136").Value = -1293.8334

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 950
$ws.Range("J94").Value = 1900
$ws.Range("L94").Value = 1900
$ws.Range("N94").Value = -2802
$ws.Range("H96").Value = 13333.333
$ws.Range("I96").Value = 13333.333
$ws.Range("K96").Value = 13333.333
$ws.Range("M96").Value = -10587.333

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5895.4287
$ws.Range("I58").Value = 3249
$ws.Range("K58").Value = 3249
$ws.Range("M58").Value = -3046
$ws.Range("H132").Value = 3366.6296
$ws.Range("I132").Value = 2906.682
$ws.Range("J132").Value = 5390.4
$ws.Range("K132").Value = 8720.045999999998
$ws.Range("L132").Value = 16171.2
$ws.Range("M132").Value = -6190.045999999998
$ws.Range("N132").Value = -21231.2
$ws.Range("H134").Value = 4023.65
$ws.Range("I134").Value = 2597
$ws.Range("J134").Value = 7352.5
$ws.Range("K134").Value = 7791
$ws.Range("L134").Value = 22057.5
$ws.Range("M134").Value = -5256
$ws.Range("N134").Value = -27127.5
$ws.Range("H136").Value = 5895.4287
$ws.Range("I136").Value = 3249
$ws.Range("K136").Value = 9747
$ws.Range("M136").Value = -7197

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1375
$ws.Range("I18").Value = 1375
$ws.Range("K18").Value = 4125
$ws.Range("M18").Value = -3956
$ws.Range("H32").Value = 14498.75
$ws.Range("J32").Value = 48999
$ws.Range("L32").Value = 146997
$ws.Range("N32").Value = -147563
$ws.Range("H58").Value = 1000
$ws.Range("J58").Value = 1000
$ws.Range("L58").Value = 3000
$ws.Range("N58").Value = -3256
$ws.Range("H97").Value = 232.66667
$ws.Range("I97").Value = 99
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 297
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = 199
$ws.Range("N97").Value = -2492
$ws.Range("H107").Value = 1423.75
$ws.Range("I107").Value = 1201
$ws.Range("J107").Value = 1582.8572
$ws.Range("K107").Value = 3603
$ws.Range("L107").Value = 4748.571599999999
$ws.Range("M107").Value = -1683
$ws.Range("N107").Value = -8588.571599999999
$ws.Range("H132").Value = 2419.9285
$ws.Range("I132").Value = 1067.6154
$ws.Range("J132").Value = 20000
$ws.Range("K132").Value = 9608.5386
$ws.Range("L132").Value = 180000
$ws.Range("M132").Value = -7078.5386
$ws.Range("N132").Value = -185060

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 44426
$ws.Range("J20").Value = 44426
$ws.Range("L20").Value = 44426
$ws.Range("N20").Value = -44916
$ws.Range("H97").Value = 1939.5927
$ws.Range("J97").Value = 1673.8572
$ws.Range("L97").Value = 1673.8572
$ws.Range("N97").Value = -2665.8572
$ws.Range("H100").Value = 44997
$ws.Range("J100").Value = 44997
$ws.Range("L100").Value = 44997
$ws.Range("N100").Value = -47161

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4008.3
$ws.Range("I7").Value = 3370.4285
$ws.Range("K7").Value = 3370.4285
$ws.Range("M7").Value = -3258.4285
$ws.Range("H22").Value = 889.4666999999999
$ws.Range("J22").Value = 848.2
$ws.Range("L22").Value = 848.2
$ws.Range("N22").Value = -1438.2
$ws.Range("H27").Value = 889.4666999999999
$ws.Range("J27").Value = 848.2
$ws.Range("L27").Value = 848.2
$ws.Range("N27").Value = -1062.2
$ws.Range("H46").Value = 2490.6875
$ws.Range("I46").Value = 985.1
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 985.1
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -797.1
$ws.Range("N46").Value = -5376
$ws.Range("H106").Value = 14166.333
$ws.Range("J106").Value = 14166.333
$ws.Range("L106").Value = 14166.333
$ws.Range("N106").Value = -16690.333
$ws.Range("H126").Value = 4008.3
$ws.Range("I126").Value = 3370.4285
$ws.Range("K126").Value = 10111.2855
$ws.Range("M126").Value = -7641.2855
$ws.Range("H132").Value = 4431.524
$ws.Range("I132").Value = 3171.8333
$ws.Range("J132").Value = 6111.1113
$ws.Range("K132").Value = 9515.499899999999
$ws.Range("L132").Value = 18333.3339
$ws.Range("M132").Value = -6985.499899999999
$ws.Range("N132").Value = -23393.3339

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 88333.336
$ws.Range("J92").Value = 88333.336
$ws.Range("L92").Value = 88333.336
$ws.Range("N92").Value = -93325.336
$ws.Range("H95").Value = 44739.75
$ws.Range("J95").Value = 44739.75
$ws.Range("L95").Value = 44739.75
$ws.Range("N95").Value = -50231.75
$ws.Range("H96").Value = 698
$ws.Range("I96").Value = 861
$ws.Range("J96").Value = 535
$ws.Range("K96").Value = 861
$ws.Range("L96").Value = 535
$ws.Range("M96").Value = 512
$ws.Range("N96").Value = -3281
$ws.Range("H97").Value = 45572
$ws.Range("J97").Value = 45572
$ws.Range("L97").Value = 45572
$ws.Range("N97").Value = -47554
